# Apply edits described in commit "add test cases for careplan":
#  - extend the "specialist" sheet (sheet3) with new columns (O:BB) on the
#    header/first data row, and three brand new data rows (3-5)
#  - add a new "ithaca" worksheet at the end of the workbook with its own
#    header row + one data row (ithacaCarePlanEligibility)
$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("specialist")

# --- specialist: new header cells on row 1 (O1:BB1) ---
$ws3.Range("O1").Value = 'Age of Menarche'
$ws3.Range("P1").Value = 'Menopausal Status'
$ws3.Range("Q1").Value = 'Number of Breast Biopsies'
$ws3.Range("R1").Value = 'Breast Density'
$ws3.Range("S1").Value = 'Number of Colon Polyps'
$ws3.Range("T1").Value = 'ReportName'
$ws3.Range("U1").Value = 'Has Had Genetic Testing'
$ws3.Range("V1").Value = 'Add Pathogenic Mutations1'
$ws3.Range("W1").Value = 'Add Pathogenic Mutations2'
$ws3.Range("X1").Value = 'Breast ER Status'
$ws3.Range("Y1").Value = 'Care Plan'
$ws3.Range("Z1").Value = 'Comorbidity'
$ws3.Range("AA1").Value = 'Age of Diagnosis'
$ws3.Range("AB1").Value = 'Race'
$ws3.Range("AC1").Value = 'Sub Race'
$ws3.Range("AD1").Value = 'Surgery'
$ws3.Range("AE1").Value = 'Gravida'
$ws3.Range("AF1").Value = 'Parity'
$ws3.Range("AG1").Value = 'Age of First Birth'
$ws3.Range("AH1").Value = 'Weight(lb)'
$ws3.Range("AI1").Value = 'Height'
$ws3.Range("AJ1").Value = 'Vendors'
$ws3.Range("AK1").Value = 'Test'
$ws3.Range("AL1").Value = 'TestName1'
$ws3.Range("AM1").Value = 'TestResultOption1'
$ws3.Range("AN1").Value = 'TestName2'
$ws3.Range("AO1").Value = 'TestResultOption2'
$ws3.Range("AP1").Value = 'FormNames'
$ws3.Range("AQ1").Value = 'TestName3'
$ws3.Range("AR1").Value = 'TestResultOption3'
$ws3.Range("AS1").Value = 'AdditionalTestResultOption'
$ws3.Range("AT1").Value = 'Family Grouping1'
$ws3.Range("AU1").Value = 'Relationship1'
$ws3.Range("AV1").Value = 'Cancer Type1'
$ws3.Range("AW1").Value = 'Family Grouping1'
$ws3.Range("AX1").Value = 'Age Of Diagnosis1'
$ws3.Range("AY1").Value = 'Relationship2'
$ws3.Range("AZ1").Value = 'Cancer Type2'
$ws3.Range("BA1").Value = 'Family Grouping2'
$ws3.Range("BB1").Value = 'Age Of Diagnosis2'

# --- specialist: new row 3 (SpecialistMamogramCheck) ---
$ws3.Range("A3").Value = 'SpecialistMamogramCheck'
$ws3.Range("B3").Value = 'YJF151102'
$ws3.Range("C3").Value = 'YJL151102'
$ws3.Range("D3").Value = 'Female'
$ws3.Range("E3").Value = '''06/08/1978'
$ws3.Range("F3").Value = 'Chicago - New User Location'
$ws3.Range("G3").Value = 'user@email.com'
$ws3.Range("H3").Value = 'user@email.com'
$ws3.Range("I3").Value = 'No'
$ws3.Range("O3").Value = 11
$ws3.Range("P3").Value = 'postmenopausal'
$ws3.Range("Q3").Value = 2
$ws3.Range("R3").Value = 'N/A'
$ws3.Range("S3").Value = 0
$ws3.Range("T3").Value = 'Patient Letter (New)'
$ws3.Range("Y3").Value = 'Mammogram (Annual) - RECOMMENDED'

# --- specialist: new row 4 (SpecialistAromataseCheck) ---
$ws3.Range("A4").Value = 'SpecialistAromataseCheck'
$ws3.Range("B4").Value = 'YJF151102'
$ws3.Range("C4").Value = 'YJL151102'
$ws3.Range("D4").Value = 'Female'
$ws3.Range("E4").Value = '''06/08/1988'
$ws3.Range("F4").Value = 'Chicago - New User Location'
$ws3.Range("G4").Value = 'user@email.com'
$ws3.Range("H4").Value = 'user@email.com'
$ws3.Range("I4").Value = 'No'
$ws3.Range("O4").Value = 11
$ws3.Range("P4").Value = 'postmenopausal'
$ws3.Range("Q4").Value = 2
$ws3.Range("R4").Value = 'N/A'
$ws3.Range("S4").Value = 0
$ws3.Range("T4").Value = 'Patient Letter (New)'
$ws3.Range("U4").Value = 'Yes'
$ws3.Range("V4").Value = 'BRCA1'
$ws3.Range("W4").Value = 'BRCA2'
$ws3.Range("X4").Value = 'Positive'
$ws3.Range("Y4").Value = 'Aromatase Inhibitors (Anastrozole 1 mg/day or Exemestane 25mg/d for 5 years) - RECOMMENDED'

# --- specialist: new row 5 (SpecialistHBOCAndHCCCheck) ---
$ws3.Range("A5").Value = 'SpecialistHBOCAndHCCCheck'
$ws3.Range("C5").Value = 'YJL151102'
$ws3.Range("D5").Value = 'Female'
$ws3.Range("E5").Value = '''06/08/1983'
$ws3.Range("F5").Value = 'Chicago - New User Location'
$ws3.Range("G5").Value = 'user@email.com'
$ws3.Range("H5").Value = 'user@email.com'
$ws3.Range("I5").Value = 'No'
$ws3.Range("P5").Value = 'premenopausal'
$ws3.Range("U5").Value = 'No'
$ws3.Range("Z5").Value = 'Diabetes Type 2 (CKD)'
$ws3.Range("AA5").Value = 39
$ws3.Range("AB5").Value = 'Asian'
$ws3.Range("AC5").Value = 'Japanese'
$ws3.Range("AD5").Value = 'Tubal Ligation'
$ws3.Range("AE5").Value = 1
$ws3.Range("AF5").Value = 1
$ws3.Range("AG5").Value = 25
$ws3.Range("AH5").Value = 178
$ws3.Range("AI5").Value = '5 ft 2 in'
$ws3.Range("AJ5").Value = 'Invitae'
$ws3.Range("AK5").Value = 'BRCA 1/2 + Multi-Cancer Panel'
$ws3.Range("AL5").Value = 'BRCA1'
$ws3.Range("AM5").Value = 'Pathogenic/Deleterious'
$ws3.Range("AN5").Value = 'BRCA2'
$ws3.Range("AO5").Value = 'VUS - Likely Pathogenic'
$ws3.Range("AP5").Value = 'UHC LMN - HBOC'
$ws3.Range("AQ5").Value = 'ATM'
$ws3.Range("AR5").Value = 'VUS - Unknown'
$ws3.Range("AS5").Value = 'VUS - Likely Benign'
$ws3.Range("AT5").Value = 'Maternal'
$ws3.Range("AU5").Value = 'Aunt'
$ws3.Range("AV5").Value = 'Ovarian'
$ws3.Range("AW5").Value = 'Maternal'
$ws3.Range("AX5").Value = 45
$ws3.Range("AY5").Value = 'Brother'
$ws3.Range("AZ5").Value = 'Colorectal'
$ws3.Range("BA5").Value = 'Immediate'
$ws3.Range("BB5").Value = 45

# --- add the new "ithaca" worksheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "ithaca"

# --- ithaca: header row ---
$ws4.Range("A1").Value = 'testCaseName'
$ws4.Range("B1").Value = 'First Name'
$ws4.Range("C1").Value = 'Last Name'
$ws4.Range("D1").Value = 'Sex'
$ws4.Range("E1").Value = 'DayOfBirth'
$ws4.Range("F1").Value = 'MonthOfBirth'
$ws4.Range("G1").Value = 'YearOfBirth'
$ws4.Range("H1").Value = 'Currently Smokes'
$ws4.Range("I1").Value = 'Years Smoked'
$ws4.Range("J1").Value = 'Number of Packs Per Day'
$ws4.Range("K1").Value = 'eligibilities'
$ws4.Range("L1").Value = 'Relationship1'
$ws4.Range("M1").Value = 'Cancer Type1'
$ws4.Range("N1").Value = 'Family Grouping1'
$ws4.Range("O1").Value = 'Age Of Diagnosis1'
$ws4.Range("P1").Value = 'Relationship2'
$ws4.Range("Q1").Value = 'Cancer Type2'
$ws4.Range("R1").Value = 'Family Grouping2'
$ws4.Range("S1").Value = 'Age Of Diagnosis2'

# --- ithaca: data row (ithacaCarePlanEligibility) ---
$ws4.Range("A2").Value = 'ithacaCarePlanEligibility'
$ws4.Range("B2").Value = 'john'
$ws4.Range("C2").Value = 'doe'
$ws4.Range("D2").Value = 'FEMALE'
$ws4.Range("E2").Value = 12
$ws4.Range("F2").Value = 12
$ws4.Range("G2").Value = 1970
$ws4.Range("H2").Value = 'Yes'
$ws4.Range("I2").Value = 11
$ws4.Range("J2").Value = 2
$ws4.Range("K2").Value = 'LUNG SCREENING'
$ws4.Range("L2").Value = 'Aunt'
$ws4.Range("M2").Value = 'Ovarian'
$ws4.Range("N2").Value = 'Maternal'
$ws4.Range("O2").Value = 45
$ws4.Range("P2").Value = 'Brother'
$ws4.Range("Q2").Value = 'Colorectal'
$ws4.Range("R2").Value = 'Immediate'
$ws4.Range("S2").Value = 45

# trailing styled-but-empty cell at T1 (keeps the sheet dimension at A1:T2,
# matching the source workbook)
$ws4.Range("T1").Font.Name = "Arial"

# the edited workbook re-opens with "specialist" as the active tab
$ws3.Activate()
